$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("report")

# --- Insert a new column before column C (this shifts C: onward right by 1) ---
$ws.Columns("C:C").Insert()

# The "value" header label in row 1 tracked the data (it was over the first data
# column). After the insert it auto-shifted from E1 to F1; move it back onto the
# new first data column C1 and clear its old shifted spot.
$headerText = $ws.Range("F1").Value2
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value2 = $headerText
$ws.Range("F1").Clear()

# --- Populate the new column C with the new "Uralsibins" / "SAN_Assessment" record ---
$ws.Range("C2").Value2 = "Uralsibins"

$ws.Range("D3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value2 = 44484

$ws.Range("C4").Value2 = "SAN_Assessment"
$ws.Range("C5").Value2 = "C:\Users\vlasenko\Documents\01.CUSTOMERS\Uralsibins\SAN Assessment\OCT2021"
$ws.Range("C6").Value2 = "C:\Users\vlasenko\Documents\06.CONFIGS\Uralsib\OCT21\FCSWITCH"
$ws.Range("C7").Value2 = "C:\Users\vlasenko\Documents\06.CONFIGS\Uralsib\OCT21\blade"

# --- sheetView selection on "report" ---
$ws.Range("C14").Select()

# --- service_tables sheet tweaks ---
$svc = $wb.Worksheets.Item("service_tables")
$svc.Range("H84").Value2 = 0
$svc.Range("F94").Value2 = "porterr_enc_crc_bad_os"
$svc.Range("G129").Value2 = 1
$svc.Range("H129").Value2 = 1

$svc.Range("H80:H81").Select()
